$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores every value as literal text (e.g. "68.440.53" or
# "1.00"), including values that look like plain numbers. Assigning such a
# string straight to .Value would let Excel auto-coerce it into a Number
# cell and silently drop formatting such as trailing zeros, so for every
# Price cell whose new value is numeric-looking we momentarily force the
# "Text" number format, assign the value, then restore the default "Normal"
# cell style so no stray formatting is left behind. (NumberFormat is set one
# cell at a time -- a multi-area Range("D1,D2").NumberFormat only reaches the
# first area on this host.)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "68.283.26"
$ws.Range("E2").Value = "  +0.76%  "

# Row 3
$ws.Range("D3").Value = "3.798.15"
$ws.Range("E3").Value = "  -0.23%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "607.30"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "163.45"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("D6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "3.802.27"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").Value = "0.516"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("E10").Value = "  -0.01%  "

# Row 11
$ws.Range("D11").Value = "6.94"
$ws.Range("E11").Value = "  +10.12%  "
$ws.Range("D11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "0.450"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("D13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "35.03"
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "4.446.29"
$ws.Range("E15").Value = "  +0.03%  "

# Row 16
$ws.Range("D16").Value = "3.819.79"
$ws.Range("E16").Value = "  +0.47%  "

# Row 17
$ws.Range("D17").Value = "68.362.41"
$ws.Range("E17").Value = "  +0.86%  "

# Row 18
$ws.Range("D18").Value = "18.11"
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("D18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "0.114"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "7.07"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "463.14"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "9.61"
$ws.Range("E22").Value = "  -2.48%  "
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "0.698"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "0.0000147"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "83.46"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "12.00"
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "2.10"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D27").Style = "Normal"

# Row 28
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D28").Style = "Normal"

# Row 29
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "9.98"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "3.957.29"
$ws.Range("E30").Value = "  +0.03%  "

# Row 31
$ws.Range("D31").Value = "2.63"
$ws.Range("E31").Value = "  -5.61%  "
$ws.Range("D31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "2.22"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "7.24"
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("D33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "29.07"
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("D34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "9.02"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = "  +1.27%  "

# Row 38
$ws.Range("D38").Value = "0.148"
$ws.Range("E38").Value = "  +7.04%  "
$ws.Range("D38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "5.88"
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("D39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "3.21"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "0.979"
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("D41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = "  +0.21%  "

# Row 43
$ws.Range("E43").Value = "  +0.00%  "

# Row 44
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "46.99"
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("D44").Style = "Normal"

# Row 45
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "152.83"
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "0.296"
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "42.88"
$ws.Range("E47").Value = "  -4.34%  "
$ws.Range("D47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "1.39"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("D48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "8.38"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "1.86"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "26.21"
$ws.Range("E51").Value = "  -6.65%  "
$ws.Range("D51").Style = "Normal"
